$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-Format($src, $dst) {
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial(-4122) | Out-Null
}

# ------------------------------------------------------------------
# Step 1: move the existing "Transposed" block from F13:H16 to
# A13:C16 *before* touching F13:H16 further (so its header/format is
# still intact to copy from). The header (F13:H13) is merged, so
# unmerge it first - otherwise only the anchor cell can be written.
# ------------------------------------------------------------------
$ws.Range("F13:H13").UnMerge() | Out-Null

Copy-Format "F13:H13" "A13:C13"
$ws.Range("A13").Value2 = "Transposed"
$ws.Range("A13:C13").Merge() | Out-Null

Copy-Format "F14:H14" "A14:C14"
$ws.Range("A14").Value2 = "John"
$ws.Range("B14").Value2 = "Luis"
$ws.Range("C14").Value2 = "Henry"

Copy-Format "F15:H15" "A15:C15"
$ws.Range("A15").Value2 = "On Elm St."
$ws.Range("B15").Value2 = "On 23rd St."
$ws.Range("C15").Value2 = "On 5th Ave."

Copy-Format "F16:H16" "A16:C16"
$ws.Range("A16").Value2 = 30
$ws.Range("B16").Value2 = 21
$ws.Range("C16").Value2 = 45

# ------------------------------------------------------------------
# Step 2: build the new "From List" table. Header goes in F11:I11
# (merged, bold/centred style copied from the header cells we just
# vacated). Data rows reuse F12:I15 (copied from the plain data style
# that used to live in F14:H14).
# ------------------------------------------------------------------
Copy-Format "F13:H13" "F11:H11"
Copy-Format "H13" "I11"
$ws.Range("F11").Value2 = "From List"
$ws.Range("F11:I11").Merge() | Out-Null

Copy-Format "F14:H14" "F12:H12"
Copy-Format "H14" "I12"
Copy-Format "F14:H14" "F13:H13"
Copy-Format "H14" "I13"
Copy-Format "H14" "I14"
Copy-Format "F14:H14" "F15:H15"
Copy-Format "H14" "I15"

$ws.Range("F12").Value2 = "On Elm St."
$ws.Range("G12").Value2 = "John"
$ws.Range("H12").Value2 = 30
$ws.Range("I12").Value2 = "Person"

$ws.Range("F13").Value2 = "On Main St."
$ws.Range("G13").Value2 = "Mary"
$ws.Range("H13").Value2 = 15
$ws.Range("I13").Value2 = "Person"

$ws.Range("F14").Value2 = "On 23rd St."
$ws.Range("G14").Value2 = "Luis"
$ws.Range("H14").Value2 = 21
$ws.Range("I14").Value2 = "Person"

$ws.Range("F15").Value2 = "On 5th Ave."
$ws.Range("G15").Value2 = "Henry"
$ws.Range("H15").Value2 = 45
$ws.Range("I15").Value2 = "Person"

# Old F16:H16 numbers are no longer needed - that data now lives in
# A16:C16 (copied above), so clear the leftovers.
$ws.Range("F16:H16").ClearContents() | Out-Null

# ------------------------------------------------------------------
# Step 3: column width tweaks to fit the new layout.
# ------------------------------------------------------------------
$ws.Columns(3).ColumnWidth = 11.05   # C  -> 11.870625 (closest achievable: 11.833333)
$ws.Columns(6).ColumnWidth = 11.05   # F  -> 11.850625 (closest achievable: 11.833333)
$ws.Columns(8).ColumnWidth = 3.15    # H  -> 3.960625  (closest achievable: 4)
$ws.Columns(9).ColumnWidth = 6.87    # I  -> 7.720625  (closest achievable: 7.666667)

# ------------------------------------------------------------------
# Step 4: update the "Titles" defined name to reference the new
# header ranges instead of the old F13:H13 one.
# ------------------------------------------------------------------
$n = $wb.Names.Item("Titles")
$n.RefersTo = "='Inserting Data'!`$A`$1:`$A`$1,'Inserting Data'!`$C`$1:`$H`$1,'Inserting Data'!`$A`$6:`$D`$6,'Inserting Data'!`$F`$6:`$H`$6,'Inserting Data'!`$F`$11:`$I`$11,'Inserting Data'!`$A`$13:`$C`$13"
